# Corrected inconsistent sample naming in test data
#
# Three samples were originally labelled "Rf-cows-pooled-d0-Rx" (lower-case
# "f" in "Rf"), inconsistent with their group label "RF-Cows pooled"
# (upper-case "RF"). Correct the sample names in column A to use the same
# "RF-Cows-pooled-d0-Rx" casing/hyphenation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A66").Value  = "RF-Cows-pooled-d0-R3"
$ws.Range("A98").Value  = "RF-Cows-pooled-d0-R2"
$ws.Range("A102").Value = "RF-Cows-pooled-d0-R1"

# Leave the sheet scrolled to the last edited row, with B102 selected,
# matching the editor's on-screen state after making the change.
$excel.ActiveWindow.ScrollRow = 79
$ws.Range("B102").Select()
